# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# Update OFF sheet (Road "R" row, row 3) with new target depth data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 421
$wsOff.Range("C3").Value = 289
$wsOff.Range("D3").Value = 74
$wsOff.Range("E3").Value = 38

# Update DEF sheet (Road "R" row, row 3) with new target depth data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 543
$wsDef.Range("C3").Value = 389
$wsDef.Range("D3").Value = 124
$wsDef.Range("E3").Value = 59
